# referenceDocsMetrics.xlsx — "Generate Errors.txt BEFORE applying special rules."
#
# Row 5 was previously computed/stamped AFTER the "special rules" pass, which
# gave it a distinct custom row style (s="3") and slightly different values
# for the "Classes without inherits from" (G) and "Array of" (L) columns.
# Re-generating Errors.txt before the special rules run changes those two
# counts and leaves row 5 formatted the same as the other plain data rows
# (2-4) instead of carrying the special, customFormat row style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two values that differ once Errors.txt is generated earlier ---
$ws.Range("G5").Value = 202
$ws.Range("L5").Value = 225

# --- Strip row 5's custom formatting so it matches the other data rows ---
# ClearFormats() wipes every cell's style in the row (including A5's date
# format), so restore A5's original date formatting by copying the format
# from A2 (which already carries the correct, pre-existing date style).
$row5 = $ws.Rows("5:5")
$row5.ClearFormats()
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Select the entire (now re-generated) row 5, matching the saved view ---
$ws.Rows("5:5").Select()
